# aggiornamento fino a 9 agosto 2021
# Append new daily rows (329-343) to Sheet1, extending the data through
# 2021-08-09 (Excel serial date 44417), matching the formatting of the
# preceding data row (row 328).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Propagate the formatting (style/number format/borders/alignment) of the
# last existing data row onto the new rows before writing values into them.
$ws.Range("A328:D328").Copy()
$ws.Range("A329:D343").PasteSpecial(-4122)

$newData = @(
    @(44403, 3, 7, 40.85205719288007),
    @(44404, 1, 8, 46.68806536329151),
    @(44405, 2, 10, 58.36008170411438),
    @(44406, 3, 12, 70.03209804493727),
    @(44407, 4, 15, 87.54012255617158),
    @(44408, 5, 19, 110.8841552378173),
    @(44409, 0, 18, 105.0481470674059),
    @(44410, 3, 18, 105.0481470674059),
    @(44411, 0, 17, 99.21213889699445),
    @(44412, 1, 16, 93.37613072658301),
    @(44413, 2, 15, 87.54012255617158),
    @(44414, 4, 15, 87.54012255617158),
    @(44415, 4, 14, 81.70411438576014),
    @(44416, 5, 19, 110.8841552378173),
    @(44417, 5, 21, 122.5561715786402)
)

$r = 329
foreach ($row in $newData) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $r = $r + 1
}
